$d = $word.ActiveDocument
$paras = $d.Paragraphs
$i = 0
$target = $null
foreach ($p in $paras) {
    if ($i -eq 14) { $target = $p }
    $i = $i + 1
}
$pEnd = $target.Range.End
Write-Output "p14 end=$pEnd text=[$($target.Range.Text)]"
$r = $d.Range($pEnd - 1, $pEnd)
Write-Output "char: $([int][char]$r.Text[0])"
$r.Delete()
Write-Output "deleted"
$r2 = $d.Range(440, 560)
Write-Output "[$($r2.Text)]"
